# Update "Contenu du stage" statistics (rows 16-20) on the Worksheet sheet.
# This reflects the new distribution of internship subjects and their
# recomputed percentages (amélioration de l'affichage des statistiques).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# C# : 15 -> 3 soutenances, 100% -> 20%
$ws.Range("E16").Value = 3
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "20 %"
$ws.Range("G16").ClearFormats()

# COBOL : 0 -> 8 soutenances, 0% -> 53.33%
$ws.Range("E17").Value = 8
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "53.33 %"
$ws.Range("G17").ClearFormats()

# ASSEMBLEUR : 0 -> 2 soutenances, 0% -> 13.33%
$ws.Range("E19").Value = 2
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "13.33 %"
$ws.Range("G19").ClearFormats()

# ANDROID : 0 -> 2 soutenances, 0% -> 13.33%
$ws.Range("E20").Value = 2
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "13.33 %"
$ws.Range("G20").ClearFormats()
